# convert µmoles to nmoles
# Append three more replicate standard-curve blocks (conc/abs pairs) below
# the existing one, reusing the same banded fill styling, and introduce two
# new shades for the "250" and one "62.5" rows that need slightly different
# fills in the new blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replicate the original A2:B8 block (with its formatting) three times ---
$src = $ws.Range("A2:B8")
$src.Copy($ws.Range("A9"))
$src.Copy($ws.Range("A16"))
$src.Copy($ws.Range("A23"))

# --- Block 1 (rows 9-15) ---
$ws.Range("A9").Value  = 500
$ws.Range("B9").Value  = 3.26
$ws.Range("A10").Value = 250
$ws.Range("B10").Value = 1.687
$ws.Range("A11").Value = 125
$ws.Range("B11").Value = 0.982
$ws.Range("A12").Value = 62.5
$ws.Range("B12").Value = 0.548
$ws.Range("A13").Value = 31.25
$ws.Range("B13").Value = 0.411
$ws.Range("A14").Value = 15.63
$ws.Range("B14").Value = 0.325
$ws.Range("A15").Value = 0
$ws.Range("B15").Value = 0.208

# --- Block 2 (rows 16-22) ---
$ws.Range("A16").Value = 500
$ws.Range("B16").Value = 3.217
$ws.Range("A17").Value = 250
$ws.Range("B17").Value = 1.687
$ws.Range("A18").Value = 125
$ws.Range("B18").Value = 0.983
$ws.Range("A19").Value = 62.5
$ws.Range("B19").Value = 0.562
$ws.Range("A20").Value = 31.25
$ws.Range("B20").Value = 0.402
$ws.Range("A21").Value = 15.63
$ws.Range("B21").Value = 0.306
$ws.Range("A22").Value = 0
$ws.Range("B22").Value = 0.181

# --- Block 3 (rows 23-29) ---
$ws.Range("A23").Value = 500
$ws.Range("B23").Value = 3.163
$ws.Range("A24").Value = 250
$ws.Range("B24").Value = 1.664
$ws.Range("A25").Value = 125
$ws.Range("B25").Value = 0.95
$ws.Range("A26").Value = 62.5
$ws.Range("B26").Value = 0.545
$ws.Range("A27").Value = 31.25
$ws.Range("B27").Value = 0.396
$ws.Range("A28").Value = 15.63
$ws.Range("B28").Value = 0.292
$ws.Range("A29").Value = 0
$ws.Range("B29").Value = 0.145

# --- New fill shades for the "250" rows in the new blocks (B10/B17/B24) ---
# fgColor FF8DBCE0 -> RGB(141,188,224) encoded as BGR-ordered OLE color
$ws.Range("B10").Interior.Color = 14728333
$ws.Range("B10").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B24").PasteSpecial(-4122)

# --- New fill shade for the "62.5" row of block 2 only (B19) ---
# fgColor FFC9E0F4 -> RGB(201,224,244) encoded as BGR-ordered OLE color
$ws.Range("B19").Interior.Color = 16048329

$excel.CutCopyMode = 0

# --- Update the active selection to match the final state ---
$ws.Range("B23:B29").Select()
